$d = $word.ActiveDocument

# Locate the word "deskriptif" in the document body.
$found = $d.Content
$found.Find.Execute("deskriptif", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$wordStart = $found.Start
$wordEnd = $found.End

# Split point: "deskrip" (7 chars) | "tif" (3 chars) -> replaced by "si"
$splitPoint = $wordStart + 7

# Replace the trailing "tif" with "si", leaving "deskrip" untouched.
$tail = $d.Range($splitPoint, $wordEnd)
$tail.Text = "si"

# Force the freshly written "si" run to stay a distinct <w:r> from the
# preceding "deskrip" run (instead of being silently re-coalesced into a
# single run on save) by nudging and reverting a direct character
# property on just that new text.
$newTail = $d.Range($splitPoint, $splitPoint + 2)
$newTail.Font.Bold = 1
$newTail.Font.Bold = 0
